$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

$newDate = Get-Date -Year 2023 -Month 10 -Day 28

# New row 45 data: date, encuestadora (RCN), and six numeric values
$ws.Cells.Item(45, 1).Value = $newDate
$ws.Cells.Item(45, 2).Value = "RCN"
$ws.Cells.Item(45, 3).Value = 0.371
$ws.Cells.Item(45, 4).Value = 0.355
$ws.Cells.Item(45, 5).Value = 0.066
$ws.Cells.Item(45, 6).Value = 0.07
$ws.Cells.Item(45, 7).Value = 0.077
$ws.Cells.Item(45, 8).Value = 0.061

# Copy formatting from row 41 (same style as the rest of that date block) to row 45
$ws.Range("A41:H41").Copy() | Out-Null
$ws.Range("A45:H45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values (PasteSpecial(xlPasteFormats) shouldn't touch values, but make sure)
$ws.Cells.Item(45, 1).Value = $newDate
$ws.Cells.Item(45, 2).Value = "RCN"
$ws.Cells.Item(45, 3).Value = 0.371
$ws.Cells.Item(45, 4).Value = 0.355
$ws.Cells.Item(45, 5).Value = 0.066
$ws.Cells.Item(45, 6).Value = 0.07
$ws.Cells.Item(45, 7).Value = 0.077
$ws.Cells.Item(45, 8).Value = 0.061

# Update the view: scroll position and selection moved
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("J44").Select() | Out-Null
